$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style from H1 so the new header cells (I1, J1)
# match the formatting (bold, border, centered) used by the other headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: row number, I value, J value
$data = @(
    @(2, 8, 8),
    @(3, 7, 7),
    @(4, 8, 9),
    @(5, 8, 8),
    @(6, 7, 7),
    @(7, 9, 9),
    @(8, 8, 8),
    @(9, 8, 8),
    @(10, 7, 7),
    @(11, 9, 9),
    @(12, 10, 10),
    @(13, 8, 8),
    @(14, 10, 10),
    @(15, 6, 7),
    @(16, 7, 7),
    @(17, 8, 8),
    @(18, 9, 9),
    @(19, 9, 10),
    @(20, 7, 7),
    @(21, 8, 8),
    @(22, 8, 8),
    @(23, 8, 8),
    @(24, 8, 8),
    @(25, 8, 8),
    @(26, 9, 9),
    @(27, 9, 9),
    @(28, 9, 9),
    @(29, 9, 9),
    @(30, 9, 9),
    @(31, 9, 9),
    @(32, 9, 9),
    @(33, 9, 9),
    @(34, 9, 9),
    @(35, 9, 9),
    @(36, 8, 9),
    @(37, 9, 9),
    @(38, 9, 9),
    @(39, 9, 9),
    @(40, 9, 9),
    @(41, 9, 9),
    @(42, 9, 9),
    @(43, 10, 10),
    @(44, 9, 9),
    @(45, 9, 9),
    @(46, 10, 10),
    @(47, 9, 9),
    @(48, 9, 9),
    @(49, 9, 9),
    @(50, 9, 9),
    @(51, 8, 9),
    @(52, 9, 9),
    @(53, 11, 11),
    @(54, 8, 9),
    @(55, 9, 9),
    @(56, 9, 9),
    @(57, 9, 9),
    @(58, 9, 9),
    @(59, 9, 9),
    @(60, 9, 9),
    @(61, 9, 9),
    @(62, 8, 9),
    @(63, 9, 9),
    @(64, 9, 9),
    @(65, 9, 9),
    @(66, 9, 9),
    @(67, 9, 9),
    @(68, 4, 4),
    @(69, 4, 4),
    @(70, 3, 3),
    @(71, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
